# Task List For RPG Hero.xlsx - apply commit changes
#
# Summary of edits (per commit message / diff):
#  - "Tasks 02-11 to 02-18" sheet: several tasks moved to "Done", time-spent /
#    over-under figures filled in, two task descriptions rewritten/expanded
#    (Finite State Machine work, coupling reduction), and a brand new task
#    ("Update MeleeWeapon to take in a prefab and itemImage") added in place
#    of the old "Change name of enemy methods to OnXXX" row - which itself
#    moved up to replace the now-removed "Have things that are being
#    effected..." task. Total-hours formula range grows to include the new
#    row. Selection moved to C9.
#  - "TaskList" sheet: the HUD/graphics task row now wraps onto two lines
#    (row height matches the row above it).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Tasks 02-11 to 02-18
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Tasks 02-11 to 02-18")

# Row 3 - "Have IconSpawner Spawn in new icons..." finished
$ws.Range("C3").Value = 0.5
$ws.Range("D3").Value = -1
$ws.Range("F3").Value = "Done"
$ws.Range("F3").Interior.Color = 5296274   # 92D050 green - "Done"

# Row 4 - "Fix Base Enemy Collisions with weapon" finished
$ws.Range("C4").Value = 0.5
$ws.Range("D4").Value = -0.5
$ws.Range("F4").Value = "Done"
$ws.Range("F4").Interior.Color = 5296274

# Row 9 - "Fix Icons and enemies colliding with each other" finished
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = -2
$ws.Range("F9").Value = "Done"
$ws.Range("F9").Interior.Color = 5296274

# Row 10 - reworded task, now in progress with hours logged
$ws.Range("A10").Value = "Research and implement Finite state machine for enemy"
$ws.Range("B10").Value = 4
$ws.Range("C10").Value = 4
$ws.Range("E10").Value = "James"
$ws.Range("F10").Value = "In Progress"
$ws.Range("F10").Interior.Color = 65535   # FFFF00 yellow - "In Progress"

# Row 11 - reworded task, finished
$ws.Range("A11").Value = "Reduce coupling/Clean up code"
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = 2.25
$ws.Range("D11").Value = 0.25
$ws.Range("E11").Value = "James"
$ws.Range("F11").Value = "Done"
$ws.Range("F11").Interior.Color = 5296274

# Row 12 - old "Change name of enemy methods to OnXXX" task, finished
$ws.Range("A12").Value = "Change name of enemy methods to OnXXX"
$ws.Range("B12").Value = 0.25
$ws.Range("C12").Value = 0.25
$ws.Range("E12").Value = "James"
$ws.Range("F12").Value = "Done"
$ws.Range("F12").Interior.Color = 5296274

# Row 13 - brand new task, finished
$ws.Range("A13").Value = "Update MeleeWeapon to take in a prefab and itemImage"
$ws.Range("B13").Value = 0.5
$ws.Range("C13").Value = 0.5
$ws.Range("E13").Value = "James"
$ws.Range("F13").Value = "Done"
$ws.Range("F13").Interior.Color = 5296274

# Totals row - hours-assigned sum now spans through the new row 13
$ws.Range("B15").Formula = "=SUM(B2:B13)"

$ws.Range("C9").Select()

# ---------------------------------------------------------------------
# Sheet: TaskList
# ---------------------------------------------------------------------
$taskList = $wb.Worksheets.Item("TaskList")
$taskList.Rows.Item(20).RowHeight = 26.25
